$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Row 1 (table row 1) ---
Replace-Text "482÷6=" "684÷9="
Replace-Text "914÷3=" "481÷7="
Replace-Text "566÷9=" "274÷4="
Replace-Text "500÷9=" "773÷4="
Replace-Text "836÷9=" "106÷2="

# --- Row 2 (table row 5) ---
Replace-Text "238÷6=" "925÷9="
Replace-Text "542÷4=" "945÷3="
Replace-Text "275÷2=" "836÷7="
Replace-Text "793÷9=" "483÷4="
Replace-Text "634÷2=" "368÷2="

# --- Row 3 (table row 9) ---
Replace-Text "901÷5=" "896÷3="
Replace-Text "295÷9=" "418÷9="
Replace-Text "601÷5=" "341÷5="
Replace-Text "523÷5=" "303÷6="
Replace-Text "557÷8=" "270÷8="

# --- Row 4 (table row 13) ---
# Original cells: 975÷3=, 547÷2=, 604÷6=, 270÷8=, 728÷8=
# New cells:      786÷6=, 280÷5=, 975÷3=, 665÷2=, 631÷2=
# Net cell count is unchanged (5 -> 5), so rewrite every cell's text in
# place rather than trying to insert/delete individual table cells.
$t = $d.Tables.Item(1)
$t.Cell(13, 1).Range.Text = "786÷6="
$t.Cell(13, 2).Range.Text = "280÷5="
$t.Cell(13, 3).Range.Text = "975÷3="
$t.Cell(13, 4).Range.Text = "665÷2="
$t.Cell(13, 5).Range.Text = "631÷2="

# --- Row 5 (table row 17) ---
Replace-Text "579÷3=" "880÷3="
Replace-Text "909÷3=" "964÷8="
Replace-Text "586÷3=" "372÷4="
Replace-Text "306÷3=" "626÷3="
Replace-Text "494÷5=" "289÷3="
